$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.086.00'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E5').Value = '  +2.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.59'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.328'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.90%  '
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1000'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '2.080.80'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '1.828.26'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.661'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = '35.048.95'
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '239.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.77%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('E24').Value = '  +2.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('E29').Value = '  +25.97%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.39%  '
$ws.Range('B32').Value = 'EURNeutrino'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
$ws.Range('D32').Value = '3.337.87'
$ws.Range('E32').Value = '  +37.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0555'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.09%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('E36').Value = '  +7.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '93.15'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.683'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.84%  '
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.310.75'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.986'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.39%  '
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('E47').Value = '  +5.74%  '
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('D49').Value = '1.996.20'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0648'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.31%  '
